$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (border on column A, date format on column E) from row 141 down to the two new rows (142-143),
# matching the per-row style pattern already used for every data row in the sheet.
$ws.Range("A141:AC141").Copy()
$ws.Range("A142:AC143").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ROW 141
$ws.Range("A141").Value = 139
$ws.Range("B141").Value = 7749757
$ws.Range("C141").Value = "India Super League"
$ws.Range("D141").Value = "India Super League"
$ws.Range("E141").Value = 45339.35416666666
$ws.Range("F141").Value = "Mohun Bagan SG"
$ws.Range("G141").Value = "Northeast United"
$ws.Range("H141").Value = 4
$ws.Range("I141").Value = 2
$ws.Range("J141").Value = "H"
$ws.Range("K141").Value = 1.666
$ws.Range("L141").Value = 3.6
$ws.Range("M141").Value = 4.333
$ws.Range("N141").Value = 1.615
$ws.Range("O141").Value = 4.2
$ws.Range("P141").Value = 4.75
$ws.Range("Q141").Value = -0.75
$ws.Range("R141").Value = 1.825
$ws.Range("S141").Value = 2.025
$ws.Range("T141").Value = 2.75
$ws.Range("U141").Value = 1.875
$ws.Range("V141").Value = 1.975
$ws.Range("W141").Value = 0.615
$ws.Range("X141").Value = -1
$ws.Range("Y141").Value = -1
$ws.Range("Z141").Value = 0.825
$ws.Range("AA141").Value = -1
$ws.Range("AB141").Value = 0.875
$ws.Range("AC141").Value = -1

# ROW 142
$ws.Range("A142").Value = 140
$ws.Range("B142").Value = 7751748
$ws.Range("C142").Value = "India Super League"
$ws.Range("D142").Value = "India Super League"
$ws.Range("E142").Value = 45339.45833333334
$ws.Range("F142").Value = "Hyderabad FC"
$ws.Range("G142").Value = "East Bengal Club"
$ws.Range("H142").Value = 0
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = "A"
$ws.Range("K142").Value = 10
$ws.Range("L142").Value = 5.25
$ws.Range("M142").Value = 1.285
$ws.Range("N142").Value = 6.5
$ws.Range("O142").Value = 4
$ws.Range("P142").Value = 1.533
$ws.Range("Q142").Value = 1
$ws.Range("R142").Value = 1.95
$ws.Range("S142").Value = 1.85
$ws.Range("T142").Value = 2.5
$ws.Range("U142").Value = 1.85
$ws.Range("V142").Value = 1.95
$ws.Range("W142").Value = -1
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = 0.5329999999999999
$ws.Range("Z142").Value = 0
$ws.Range("AA142").Value = 0
$ws.Range("AB142").Value = -1
$ws.Range("AC142").Value = 0.95

# ROW 143
$ws.Range("A143").Value = 141
$ws.Range("B143").Value = 7749869
$ws.Range("C143").Value = "India Super League"
$ws.Range("D143").Value = "India Super League"
$ws.Range("E143").Value = 45340.45833333334
$ws.Range("F143").Value = "Mumbai City FC"
$ws.Range("G143").Value = "Bengaluru"
$ws.Range("H143").Value = 2
$ws.Range("I143").Value = 0
$ws.Range("J143").Value = "H"
$ws.Range("K143").Value = 1.65
$ws.Range("L143").Value = 4
$ws.Range("M143").Value = 5
$ws.Range("N143").Value = 1.5
$ws.Range("O143").Value = 4.333
$ws.Range("P143").Value = 6
$ws.Range("Q143").Value = -1
$ws.Range("R143").Value = 1.85
$ws.Range("S143").Value = 2
$ws.Range("T143").Value = 2.75
$ws.Range("U143").Value = 1.925
$ws.Range("V143").Value = 1.925
$ws.Range("W143").Value = 0.5
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = 0.8500000000000001
$ws.Range("AA143").Value = -1
$ws.Range("AB143").Value = -1
$ws.Range("AC143").Value = 0.925

